# Attendance scanner app re-saved this log workbook: the recorded "Log Time"
# for the General Surgery entry (row 2) was corrected/updated.
# Apply it the way a user would in Excel: format D2 as a time (h:mm:ss),
# give it its numeric time-of-day value, and make sure the font color is
# explicit black (as the scanner app writes it), then leave that cell
# selected/active, matching the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timeCell = $ws.Range("D2")
$timeCell.NumberFormat = "h:mm:ss"
$timeCell.Value = 0.46059027777777778
$timeCell.Font.Color = 0

$timeCell.Select()
